$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous row (7) had a stray blank "MATRICULA REGISTRO CIVIL" cell
# (AW7) left over from an earlier edit; clear it out as part of this cleanup.
$ws.Range("AW7").ClearContents()

# --- Append new student registration as row 8 -----------------------------
# Columns whose values are numeric-looking or date-looking strings (IDs,
# CEP/phone with leading zeros, plain-text dates) need the cell pre-formatted
# as Text, otherwise Excel auto-converts them to numbers/dates on entry and
# the literal text (e.g. a leading "0") would be lost.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("K8").NumberFormat = "@"
$ws.Range("M8").NumberFormat = "@"
$ws.Range("N8").NumberFormat = "@"
$ws.Range("AE8").NumberFormat = "@"
$ws.Range("AG8").NumberFormat = "@"
$ws.Range("AI8").NumberFormat = "@"
$ws.Range("AN8").NumberFormat = "@"
$ws.Range("AO8").NumberFormat = "@"
$ws.Range("AP8").NumberFormat = "@"
$ws.Range("AW8").NumberFormat = "@"
$ws.Range("AX8").NumberFormat = "@"
$ws.Range("AY8").NumberFormat = "@"
$ws.Range("AZ8").NumberFormat = "@"
$ws.Range("BA8").NumberFormat = "@"

$ws.Range("A8").Value = "29129"
$ws.Range("B8").Value = "aluno teste"
$ws.Range("C8").Value = "1"
$ws.Range("D8").Value = "Branca"
$ws.Range("E8").Value = "Masculino"
$ws.Range("F8").Value = "PA"
$ws.Range("G8").Value = "Belem"
$ws.Range("H8").Value = "PA"
$ws.Range("I8").Value = "Teste"
$ws.Range("J8").Value = "Belem"
$ws.Range("K8").Value = "2000-01-05"
$ws.Range("L8").Value = "PC"
$ws.Range("M8").Value = "1"
$ws.Range("N8").Value = "11111111111"
$ws.Range("O8").Value = "SIM"
$ws.Range("P8").Value = "NÃO"
$ws.Range("Q8").Value = "NÃO"
$ws.Range("R8").Value = "NÃO"
$ws.Range("S8").Value = "SIM"
$ws.Range("T8").Value = "SIM"
$ws.Range("U8").Value = "NÃO"
$ws.Range("V8").Value = "NÃO"
$ws.Range("W8").Value = "NÃO"
$ws.Range("X8").Value = "NÃO"
$ws.Range("Y8").Value = "NÃO"
$ws.Range("Z8").Value = "NÃO"
$ws.Range("AA8").Value = "NÃO"
$ws.Range("AB8").Value = "SIM"
$ws.Range("AC8").Value = "Rua Teste"
$ws.Range("AD8").Value = "Teste 123"
$ws.Range("AE8").Value = "123"
$ws.Range("AF8").Value = "Teste"
$ws.Range("AG8").Value = "123456"
$ws.Range("AH8").Value = "Urbana"
$ws.Range("AI8").Value = "01290129"
$ws.Range("AJ8").Value = "teste@gmail.com"
$ws.Range("AK8").Value = "Teste"
$ws.Range("AL8").Value = "Teste"
$ws.Range("AM8").Value = "Teste"
$ws.Range("AN8").Value = "11111"
$ws.Range("AO8").Value = "1/2/2024"
$ws.Range("AP8").Value = "11111"
$ws.Range("AQ8").Value = "29/1/2024"
$ws.Range("AR8").Value = "Manhã"
$ws.Range("AS8").Value = "01. Berçário I"
$ws.Range("AT8").Value = "03 - Escola Estadual"
$ws.Range("AU8").Value = "SIM"
$ws.Range("AV8").Value = "SIM"
$ws.Range("AW8").Value = "1"
$ws.Range("AX8").Value = "1"
$ws.Range("AY8").Value = "1"
$ws.Range("AZ8").Value = "1"
$ws.Range("BA8").Value = "2000-01-05"
